{"js": "// Office.js (Word JavaScript API) edit script.\n// Goal (per commit message \"doc: added ID byte to the legend\"):\n//   1. Change \"...inviati 8 byte al bridge...\" -> \"...inviati 9 byte al bridge...\"\n//   2. Insert a new legend paragraph right after the intro paragraph:\n//        \"0: ID-> identificativo della cassetta (Local -> Bridge)\"\n//      (placed immediately before the existing \"1: Presenza/assenza...\" paragraph)\n\n// --- Step 1: bump the byte count from 8 to 9 in the intro sentence ---\nconst introMatches = context.document.body.search(\"inviati 8 byte\", { matchCase: true });\nintroMatches.load(\"items\");\nawait context.sync();\n\nif (introMatches.items.length > 0) {\n  introMatches.items[0].insertText(\"inviati 9 byte\", \"Replace\");\n  await context.sync();\n}\n\n// --- Step 2: insert the new \"0: ID-> ...\" legend paragraph ---\nconst anchorMatches = context.document.body.search(\n  \"Di seguito la legenda sul significato dei seguenti dati:\",\n  { matchCase: true }\n);\nanchorMatches.load(\"items\");\nawait context.sync();\n\nconst introParagraph = anchorMatches.items[0].paragraphs.getFirst();\nconst newParagraph = introParagraph.insertParagraph(\n  \"0: ID-> identificativo della cassetta (Local -> Bridge)\",\n  \"After\"\n);\n\n// Match the paragraph justification used by the rest of the legend list.\nnewParagraph.alignment = \"Justified\";\n\nawait context.sync();\n", "ps1": "# Word COM interop edit script.\n# Goal (per commit message \"doc: added ID byte to the legend\"):\n#   1. Change \"...inviati 8 byte al bridge...\" -> \"...inviati 9 byte al bridge...\"\n#   2. Insert a new legend paragraph right after the intro paragraph:\n#        \"0: ID-> identificativo della cassetta (Local -> Bridge)\"\n#      (placed immediately before the existing \"1: Presenza/assenza...\" paragraph)\n\n$d = $word.ActiveDocument\n\n# --- Step 1: bump the byte count from 8 to 9 in the intro sentence ---\n$findRange = $d.Content\n$findRange.Find.Execute(\"inviati 8 byte\", $false, $false, $false, $false, $false, $true, 1, $false, \"inviati 9 byte\", 2) | Out-Null\n\n# --- Step 2: insert the new \"0: ID-> ...\" legend paragraph right after the intro paragraph ---\n$anchorRange = $d.Content\n$anchorRange.Find.Execute(\"Di seguito la legenda sul significato dei seguenti dati:\") | Out-Null\n$introParagraph = $anchorRange.Paragraphs(1)\n$introParagraph.Range.InsertParagraphAfter()\n\n# Locate the freshly-inserted (now blank) paragraph and give it its text.\n$paragraphs = $d.Paragraphs\nfor ($i = 1; $i -le $paragraphs.Count; $i++) {\n    $paraText = $paragraphs.Item($i).Range.Text\n    if ($paraText -like \"*Di seguito la legenda sul significato dei seguenti dati:*\") {\n        $newParagraph = $paragraphs.Item($i + 1)\n        $newParagraph.Range.Text = \"0: ID-> identificativo della cassetta (Local -> Bridge)\"\n        break\n    }\n}\n"}
